$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
# ---------------------------------------------------------------------------
$srcSheet  = $wb.Worksheets.Item("2021-Q4")
$newSheet  = $wb.Worksheets.Add([System.Type]::Missing, $srcSheet)
$newSheet.Name = "2022-Q1"

# Copy header-row (B1:H1) formatting and A-column formatting from an existing
# data sheet so the new sheet matches the look of its siblings.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$data = @(
    @(0, "003993", "前海开源沪港深核心驱动灵活配置混合",        "0.58", "82.10", "8.58", "0.0498", 1),
    @(1, "010707", "安信平稳合盈一年持有期混合A",                "9.04", "20.56", "0.40", "0.0362", 9),
    @(2, "004316", "前海开源沪港深裕鑫灵活配置混合A",            "0.64", "90.55", "3.15", "0.0202", 3),
    @(3, "004317", "前海开源沪港深裕鑫灵活配置混合C",            "0.47", "90.55", "3.15", "0.0148", 3),
    @(4, "161124", "易方达香港恒生综合小型股指数（QDII-LOF）A",  "0.28", "92.62", "1.43", "0.0040", 6),
    @(5, "010708", "安信平稳合盈一年持有期混合C",                "0.26", "20.56", "0.40", "0.0010", 9),
    @(6, "006263", "易方达香港恒生综合小型股指数（QDII-LOF）C",  "0.06", "92.62", "1.43", "0.0009", 6)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Add a matching summary row to "总计"
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 7
$summary.Cells.Item(2, 4).Value = 0.13

# Re-number the index column (A) sequentially, since Insert() only shifts
# the existing static values down rather than recalculating them.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

# Restore the originally active sheet/selection (this edit should not change
# which tab is shown when the workbook is opened).
$originalActiveSheet.Activate()
$null = $originalActiveSheet.Range("A1").Select()

Write-Host "done"
